# Commit: "Add default population for worksWith"
#
# Adds a new "workswith" relation to the #People sheet: two extra columns
# (F and G) that hold, per person row, the code(s) of the person(s) they
# work with. Also mirrors the header cells' look (style) from the existing
# empty F1/F2 placeholder cells, and leaves the workbook with the #People
# sheet active/selected (matching the author's last on-screen state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#People")

# --- Headers (row 1 / row 2) ---------------------------------------------
# F1/F2 already exist as empty, styled placeholder cells; G1/G2 are brand
# new, so first clone F's look onto G before writing the values.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("F1").Value = "workswith"
$ws.Range("G1").Value = "workswith"
$ws.Range("F2").Value = "Person"
$ws.Range("G2").Value = "Person"

# --- Default "workswith" population (rows 3-5) ---------------------------
# p10001 = A. Arends (row 3), p10002 = B. Billekens (row 4), p10003 = C.
# Curly (row 5). Arends works with both Billekens and Curly; Billekens and
# Curly each work with Arends.
# These columns previously held left-over (empty) date-formatted cells, so
# clear that formatting first to get plain/general-formatted text cells.
$ws.Range("F3").ClearFormats()
$ws.Range("F4").ClearFormats()
$ws.Range("F5").ClearFormats()

$ws.Range("F3").Value = "p10002"
$ws.Range("G3").Value = "p10003"
$ws.Range("F4").Value = "p10001"
$ws.Range("F5").Value = "p10001"

# --- Window / selection state ---------------------------------------------
# The author ended up on the #People sheet with F10 selected, while the
# #Projects sheet's remembered selection moved to H2.
$ws1 = $wb.Worksheets.Item("#Projects")
$ws1.Activate() | Out-Null
$ws1.Range("H2").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("F10").Select() | Out-Null
